$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency table (coin name / link / price / 1h volume
# change) with the latest scraped values from the commit.
#
# Price cells (column D) hold numeric-looking text ("244.90", "1.726.56", ...)
# that must stay TEXT (matching the original inlineStr cells) instead of being
# auto-coerced by Excel into numbers, which would silently drop meaningful
# trailing zeros (e.g. "244.90" -> 244.9) or mis-handle the thousands-dot
# notation (e.g. "26.450.75"). Forcing the Text number format before the
# write keeps them as strings, exactly like the source data.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.450.75'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.726.56'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.90'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4792'
$ws.Range('E7').Value = '  +1.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2681'
$ws.Range('E8').Value = '  +1.33%  '
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.730.10'
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07132'
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.69'
$ws.Range('E12').Value = '  +2.72%  '
$ws.Range('E13').Value = '  +4.45%  '
$ws.Range('E14').Value = '  +2.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.20'
$ws.Range('E15').Value = '  +1.41%  '
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.463.25'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('E19').Value = '  +1.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.70'
$ws.Range('E20').Value = '  +0.67%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.952.93'
$ws.Range('E21').Value = '  +2.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.540'
$ws.Range('E22').Value = '  -0.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.921'
$ws.Range('E23').Value = '  +0.78%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.297'
$ws.Range('E24').Value = '  -0.85%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.35'
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.36'
$ws.Range('E26').Value = '  +1.22%  '
$ws.Range('E27').Value = '  +2.24%  '
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.980'
$ws.Range('E30').Value = '  -1.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08013'
$ws.Range('E31').Value = '  +3.13%  '
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04553'
$ws.Range('E33').Value = '  +2.97%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.619'
$ws.Range('E34').Value = '  +0.31%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6369'
$ws.Range('E35').Value = '  +2.31%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9886'
$ws.Range('E36').Value = '  +1.65%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9342'
$ws.Range('E37').Value = '  +1.88%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.045'
$ws.Range('E38').Value = '  +6.72%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.411'
$ws.Range('E39').Value = '  -0.49%  '
$ws.Range('B40').Value = 'Quant'
$ws.Range('C40').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '106.74'
$ws.Range('E40').Value = '  -4.12%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.003'
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01498'
$ws.Range('E42').Value = '  +1.78%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.680'
$ws.Range('E43').Value = '  +10.57%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3906'
$ws.Range('E44').Value = '  +2.41%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.992'
$ws.Range('E45').Value = '  +11.78%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1191'
$ws.Range('E46').Value = '  +4.26%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05322'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '31.06'
$ws.Range('E48').Value = '  +1.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.854'
$ws.Range('E49').Value = '  +2.38%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.268'
$ws.Range('E50').Value = '  +3.63%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3421'
$ws.Range('E51').Value = '  +1.10%  '
